$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = '41.141.86'
$ws.Cells.Item(2,5).Value = '  -2.26%  '

$ws.Cells.Item(3,4).Value = '2.180.35'
$ws.Cells.Item(3,5).Value = '  -2.05%  '

$ws.Cells.Item(4,5).Value = '  +0.16%  '

$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = '237.05'
$ws.Cells.Item(5,5).Value = '  -2.63%  '

$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = '0.615'
$ws.Cells.Item(6,5).Value = '  -2.15%  '

$ws.Cells.Item(7,4).NumberFormat = "@"
$ws.Cells.Item(7,4).Value = '71.11'
$ws.Cells.Item(7,5).Value = '  -3.33%  '

$ws.Cells.Item(8,5).Value = '  +0.01%  '

$ws.Cells.Item(9,4).NumberFormat = "@"
$ws.Cells.Item(9,4).Value = '0.580'
$ws.Cells.Item(9,5).Value = '  -5.28%  '

$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value = '40.53'
$ws.Cells.Item(10,5).Value = '  -6.19%  '

$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value = '0.0931'
$ws.Cells.Item(11,5).Value = '  -3.63%  '

$ws.Cells.Item(12,2).Value = 'TRON'
$ws.Cells.Item(12,3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Cells.Item(12,4).NumberFormat = "@"
$ws.Cells.Item(12,4).Value = '0.102'
$ws.Cells.Item(12,5).Value = '  -1.71%  '

$ws.Cells.Item(13,2).Value = 'Polkadot'
$ws.Cells.Item(13,3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(13,4).NumberFormat = "@"
$ws.Cells.Item(13,4).Value = '6.81'
$ws.Cells.Item(13,5).Value = '  -4.37%  '

$ws.Cells.Item(14,4).Value = '2.505.39'
$ws.Cells.Item(14,5).Value = '  -1.94%  '

$ws.Cells.Item(15,4).NumberFormat = "@"
$ws.Cells.Item(15,4).Value = '14.02'
$ws.Cells.Item(15,5).Value = '  -2.07%  '

$ws.Cells.Item(16,4).NumberFormat = "@"
$ws.Cells.Item(16,4).Value = '0.815'
$ws.Cells.Item(16,5).Value = '  -3.61%  '

$ws.Cells.Item(17,4).Value = '2.198.44'
$ws.Cells.Item(17,5).Value = '  -1.09%  '

$ws.Cells.Item(18,4).Value = '41.054.34'
$ws.Cells.Item(18,5).Value = '  -2.17%  '

$ws.Cells.Item(19,4).NumberFormat = "@"
$ws.Cells.Item(19,4).Value = '0.0000103'
$ws.Cells.Item(19,5).Value = '  -7.94%  '

$ws.Cells.Item(20,2).Value = 'Litecoin'
$ws.Cells.Item(20,3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(20,4).NumberFormat = "@"
$ws.Cells.Item(20,4).Value = '70.76'
$ws.Cells.Item(20,5).Value = '  -2.42%  '

$ws.Cells.Item(21,2).Value = 'Uniswap'
$ws.Cells.Item(21,3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Cells.Item(21,4).NumberFormat = "@"
$ws.Cells.Item(21,4).Value = '6.00'
$ws.Cells.Item(21,5).Value = '  -2.99%  '

$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(22,4).Value = '10.11'
$ws.Cells.Item(22,5).Value = '  -1.92%  '

$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(23,4).Value = '226.90'
$ws.Cells.Item(23,5).Value = '  -1.59%  '

$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value = '1.98'
$ws.Cells.Item(24,5).Value = '  -7.27%  '

$ws.Cells.Item(25,5).Value = '  +0.06%  '

$ws.Cells.Item(26,4).NumberFormat = "@"
$ws.Cells.Item(26,4).Value = '11.01'
$ws.Cells.Item(26,5).Value = '  -6.79%  '

$ws.Cells.Item(27,4).NumberFormat = "@"
$ws.Cells.Item(27,4).Value = '3.59'
$ws.Cells.Item(27,5).Value = '  -0.08%  '

$ws.Cells.Item(28,4).NumberFormat = "@"
$ws.Cells.Item(28,4).Value = '2.22'
$ws.Cells.Item(28,5).Value = '  -2.80%  '

$ws.Cells.Item(29,4).NumberFormat = "@"
$ws.Cells.Item(29,4).Value = '2.19'
$ws.Cells.Item(29,5).Value = '  -1.19%  '

$ws.Cells.Item(30,4).NumberFormat = "@"
$ws.Cells.Item(30,4).Value = '166.68'
$ws.Cells.Item(30,5).Value = '  +0.07%  '

$ws.Cells.Item(31,4).NumberFormat = "@"
$ws.Cells.Item(31,4).Value = '20.05'
$ws.Cells.Item(31,5).Value = '  -3.10%  '

$ws.Cells.Item(32,4).NumberFormat = "@"
$ws.Cells.Item(32,4).Value = '31.46'
$ws.Cells.Item(32,5).Value = '  +6.08%  '

$ws.Cells.Item(33,4).NumberFormat = "@"
$ws.Cells.Item(33,4).Value = '0.0784'
$ws.Cells.Item(33,5).Value = '  -1.40%  '

$ws.Cells.Item(34,4).NumberFormat = "@"
$ws.Cells.Item(34,4).Value = '5.17'
$ws.Cells.Item(34,5).Value = '  -6.72%  '

$ws.Cells.Item(35,4).NumberFormat = "@"
$ws.Cells.Item(35,4).Value = '0.121'
$ws.Cells.Item(35,5).Value = '  -2.69%  '

$ws.Cells.Item(36,4).NumberFormat = "@"
$ws.Cells.Item(36,4).Value = '0.106'
$ws.Cells.Item(36,5).Value = '  -8.88%  '

$ws.Cells.Item(37,4).NumberFormat = "@"
$ws.Cells.Item(37,4).Value = '4.14'
$ws.Cells.Item(37,5).Value = '  -5.86%  '

$ws.Cells.Item(38,4).NumberFormat = "@"
$ws.Cells.Item(38,4).Value = '0.0290'
$ws.Cells.Item(38,5).Value = '  -4.19%  '

$ws.Cells.Item(39,4).NumberFormat = "@"
$ws.Cells.Item(39,4).Value = '12.55'
$ws.Cells.Item(39,5).Value = '  -4.97%  '

$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,4).Value = '2.07'
$ws.Cells.Item(40,5).Value = '  -3.51%  '

$ws.Cells.Item(41,4).NumberFormat = "@"
$ws.Cells.Item(41,4).Value = '5.47'
$ws.Cells.Item(41,5).Value = '  -2.02%  '

$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(42,4).Value = '60.66'
$ws.Cells.Item(42,5).Value = '  -7.33%  '

$ws.Cells.Item(43,4).NumberFormat = "@"
$ws.Cells.Item(43,4).Value = '0.192'
$ws.Cells.Item(43,5).Value = '  -4.10%  '

$ws.Cells.Item(44,4).NumberFormat = "@"
$ws.Cells.Item(44,4).Value = '8.43'
$ws.Cells.Item(44,5).Value = '  -3.45%  '

$ws.Cells.Item(45,2).Value = 'Cronos'
$ws.Cells.Item(45,3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(45,4).NumberFormat = "@"
$ws.Cells.Item(45,4).Value = '0.0977'
$ws.Cells.Item(45,5).Value = '  -3.80%  '

$ws.Cells.Item(46,2).Value = 'Aave'
$ws.Cells.Item(46,3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(46,4).NumberFormat = "@"
$ws.Cells.Item(46,4).Value = '99.69'
$ws.Cells.Item(46,5).Value = '  -5.06%  '

$ws.Cells.Item(47,4).NumberFormat = "@"
$ws.Cells.Item(47,4).Value = '1.10'
$ws.Cells.Item(47,5).Value = '  -2.55%  '

$ws.Cells.Item(48,4).NumberFormat = "@"
$ws.Cells.Item(48,4).Value = '1.14'
$ws.Cells.Item(48,5).Value = '  -2.38%  '

$ws.Cells.Item(49,4).NumberFormat = "@"
$ws.Cells.Item(49,4).Value = '2.25'
$ws.Cells.Item(49,5).Value = '  -6.49%  '

$ws.Cells.Item(50,4).NumberFormat = "@"
$ws.Cells.Item(50,4).Value = '2.66'
$ws.Cells.Item(50,5).Value = '  -2.10%  '

$ws.Cells.Item(51,4).Value = '2.381.97'
$ws.Cells.Item(51,5).Value = '  -2.09%  '
